$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk 0: ALC row 11 (anchor G=5533)
$ws.Range("H11").Value = 33333416
$ws.Range("I11").Value = 33333416
$ws.Range("K11").Value = 33333416
$ws.Range("M11").Value = -33333276

# hunk 1: ALC row 29 (anchor G=4575)
$ws.Range("H29").Value = 2307.5264
$ws.Range("J29").Value = 3668.6667
$ws.Range("L29").Value = 11006.0001
$ws.Range("N29").Value = -11568.0001

# hunk 2: ALC row 51 (anchor G=5486)
$ws.Range("H51").Value = 9618.5
$ws.Range("I51").Value = 9499.75
$ws.Range("J51").Value = 9737.25
$ws.Range("K51").Value = 9499.75
$ws.Range("L51").Value = 9737.25
$ws.Range("M51").Value = -9015.75
$ws.Range("N51").Value = -10705.25

# hunk 3: ALC row 62 (anchor G=27781)
$ws.Range("H62").Value = 11175.865
$ws.Range("I62").Value = 12962.4
$ws.Range("J62").Value = 9074.058999999999
$ws.Range("K62").Value = 12962.4
$ws.Range("L62").Value = 9074.058999999999
$ws.Range("M62").Value = -12338.4
$ws.Range("N62").Value = -10322.059

# hunk 4: ALC row 64 (anchor G=5506)
$ws.Range("H64").Value = 5119
$ws.Range("J64").Value = 5148.75
$ws.Range("L64").Value = 5148.75
$ws.Range("N64").Value = -5644.75

# hunk 5: ALC row 65 (anchor G=27781)
$ws.Range("H65").Value = 11175.865
$ws.Range("I65").Value = 12962.4
$ws.Range("J65").Value = 9074.058999999999
$ws.Range("K65").Value = 64812
$ws.Range("L65").Value = 45370.295
$ws.Range("M65").Value = -61692
$ws.Range("N65").Value = -51610.295

# hunk 6: ALC row 67 (anchor G=5506)
$ws.Range("H67").Value = 5119
$ws.Range("J67").Value = 5148.75
$ws.Range("L67").Value = 5148.75
$ws.Range("N67").Value = -6864.75

# hunk 7: ALC row 116 (anchor G=27778)
$ws.Range("H116").Value = 4388.4116
$ws.Range("I116").Value = 5749.8335
$ws.Range("J116").Value = 4096.6787
$ws.Range("K116").Value = 5749.8335
$ws.Range("L116").Value = 4096.6787
$ws.Range("M116").Value = -2307.8335
$ws.Range("N116").Value = -10980.6787

# hunk 8: ALC row 125 (anchor G=36228)
$ws.Range("H125").Value = 1237.8
$ws.Range("I125").Value = 896.6667
$ws.Range("K125").Value = 8070.0003
$ws.Range("M125").Value = -5610.0003

# hunk 9: ALC row 137 (anchor G=44013)
$ws.Range("H137").Value = 10643.333
$ws.Range("I137").Value = 993
$ws.Range("K137").Value = 2979
$ws.Range("M137").Value = -429

$ws = $wb.Worksheets.Item("ARM")
# hunk 10: ARM row 17 (anchor G=2495)
$ws.Range("H17").Value = 8
$ws.Range("I17").Value = 8
$ws.Range("K17").Value = 8
$ws.Range("M17").Value = 165

# hunk 11: ARM row 110 (anchor G=27708)
$ws.Range("H110").Value = 3077.0833
$ws.Range("J110").Value = 2368.3333
$ws.Range("L110").Value = 2368.3333
$ws.Range("N110").Value = -6458.3333

# hunk 12: ARM row 139 (anchor G=42321)
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

$ws = $wb.Worksheets.Item("BSM")
# hunk 13: BSM row 99 (anchor G=19943)
$ws.Range("H99").Value = 54563.7
$ws.Range("I99").Value = 87284.414
$ws.Range("J99").Value = 5482.625
$ws.Range("K99").Value = 87284.414
$ws.Range("L99").Value = 5482.625
$ws.Range("M99").Value = -85786.414
$ws.Range("N99").Value = -8478.625

# hunk 14: BSM row 134 (anchor G=43998)
$ws.Range("H134").Value = 1799.4482
$ws.Range("I134").Value = 1710.5555
$ws.Range("K134").Value = 5131.666499999999
$ws.Range("M134").Value = -2596.666499999999

$ws = $wb.Worksheets.Item("CRP")
# hunk 15: CRP row 42 (anchor G=1847)
$ws.Range("H42").Value = 15000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 15000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 15000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -16186

# hunk 16: CRP row 86 (anchor G=12584)
$ws.Range("H86").Value = 25528.053
$ws.Range("I86").Value = 40715.273
$ws.Range("J86").Value = 4645.625
$ws.Range("K86").Value = 40715.273
$ws.Range("L86").Value = 4645.625
$ws.Range("M86").Value = -39592.273
$ws.Range("N86").Value = -6891.625

# hunk 17: CRP row 89 (anchor G=12584)
$ws.Range("H89").Value = 25528.053
$ws.Range("I89").Value = 40715.273
$ws.Range("J89").Value = 4645.625
$ws.Range("K89").Value = 203576.365
$ws.Range("L89").Value = 23228.125
$ws.Range("M89").Value = -197960.365
$ws.Range("N89").Value = -34460.125

# hunk 18: CRP row 94 (anchor G=32934)
$ws.Range("H94").Value = 3510.25
$ws.Range("I94").Value = 3474.5
$ws.Range("J94").Value = 3522.1667
$ws.Range("K94").Value = 3474.5
$ws.Range("L94").Value = 3522.1667
$ws.Range("M94").Value = -3023.5
$ws.Range("N94").Value = -4424.1667

# hunk 19: CRP row 107 (anchor G=27689)
$ws.Range("H107").Value = 487
$ws.Range("I107").Value = 418
$ws.Range("K107").Value = 418
$ws.Range("M107").Value = 1502

# hunk 20: CRP row 111 (anchor G=25792)
$ws.Range("H111").Value = 62499
$ws.Range("I111").Value = 39999
$ws.Range("J111").Value = 69999
$ws.Range("K111").Value = 39999
$ws.Range("L111").Value = 69999
$ws.Range("M111").Value = -35909
$ws.Range("N111").Value = -78179

# hunk 21: CRP row 130 (anchor G=34689)
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# hunk 22: CRP row 141 (anchor G=43345)
$ws.Range("H141").Value = 824997.25
$ws.Range("J141").Value = 824997.25
$ws.Range("L141").Value = 824997.25
$ws.Range("N141").Value = -835357.25

$ws = $wb.Worksheets.Item("CUL")
# hunk 23: CUL row 26 (anchor G=4746)
$ws.Range("H26").Value = 650
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# hunk 24: CUL row 60 (anchor G=4750)
$ws.Range("H60").Value = 246.33333
$ws.Range("I60").Value = 246.33333
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 738.99999
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -487.99999
$ws.Range("N60").ClearContents()

# hunk 25: CUL row 64 (anchor G=12861)
$ws.Range("H64").Value = 3412.25
$ws.Range("J64").Value = 4333
$ws.Range("L64").Value = 12999
$ws.Range("N64").Value = -13539

# hunk 26: CUL row 67 (anchor G=12861)
$ws.Range("H67").Value = 3412.25
$ws.Range("J67").Value = 4333
$ws.Range("L67").Value = 12999
$ws.Range("N67").Value = -14871

# hunk 27: CUL row 122 (anchor G=36078)
$ws.Range("H122").Value = 344.22726
$ws.Range("I122").Value = 319.33334
$ws.Range("J122").Value = 374.1
$ws.Range("K122").Value = 2874.00006
$ws.Range("L122").Value = 3366.9
$ws.Range("M122").Value = -424.0000600000003
$ws.Range("N122").Value = -8266.9

# hunk 28: CUL row 124 (anchor G=36040)
$ws.Range("H124").Value = 1620
$ws.Range("I124").Value = 1620
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 4860
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 50
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# hunk 29: GSM row 80 (anchor G=12521)
$ws.Range("H80").Value = 3031.2666
$ws.Range("I80").Value = 2606.1667
$ws.Range("J80").Value = 3314.6667
$ws.Range("K80").Value = 2606.1667
$ws.Range("L80").Value = 3314.6667
$ws.Range("M80").Value = -1608.1667
$ws.Range("N80").Value = -5310.6667

# hunk 30: GSM row 83 (anchor G=12521)
$ws.Range("H83").Value = 3031.2666
$ws.Range("I83").Value = 2606.1667
$ws.Range("J83").Value = 3314.6667
$ws.Range("K83").Value = 13030.8335
$ws.Range("L83").Value = 16573.3335
$ws.Range("M83").Value = -8038.833500000001
$ws.Range("N83").Value = -26557.3335

# hunk 31: GSM row 126 (anchor G=36184)
$ws.Range("H126").Value = 4921.2383
$ws.Range("J126").Value = 6899.5
$ws.Range("L126").Value = 20698.5
$ws.Range("N126").Value = -25638.5

# hunk 32: GSM row 132 (anchor G=44008)
$ws.Range("H132").Value = 36362.867
$ws.Range("I132").Value = 40659.46
$ws.Range("K132").Value = 121978.38
$ws.Range("M132").Value = -119448.38

# hunk 33: GSM row 141 (anchor G=42504)
$ws.Range("H141").Value = 69000
$ws.Range("J141").Value = 69000
$ws.Range("L141").Value = 69000
$ws.Range("N141").Value = -79360

$ws = $wb.Worksheets.Item("LTW")
# hunk 34: LTW row 55 (anchor G=5284)
$ws.Range("H55").Value = 1262.4286
$ws.Range("I55").Value = 967.7143
$ws.Range("J55").Value = 1851.8572
$ws.Range("K55").Value = 967.7143
$ws.Range("L55").Value = 1851.8572
$ws.Range("M55").Value = -794.7143
$ws.Range("N55").Value = -2197.8572

# hunk 35: LTW row 100 (anchor G=19995)
$ws.Range("H100").Value = 3474.9092
$ws.Range("I100").Value = 3298.7646
$ws.Range("K100").Value = 3298.7646
$ws.Range("M100").Value = -2757.7646

# hunk 36: LTW row 136 (anchor G=44060)
$ws.Range("H136").Value = 3285.5588
$ws.Range("I136").Value = 2766.28
$ws.Range("K136").Value = 8298.84
$ws.Range("M136").Value = -5748.84

$ws = $wb.Worksheets.Item("WVR")
# hunk 37: WVR row 19 (anchor G=2666)
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# hunk 38: WVR row 20 (anchor G=3023)
$ws.Range("H20").Value = 8000
$ws.Range("I20").Value = 6000
$ws.Range("K20").Value = 6000
$ws.Range("M20").Value = -5760

# hunk 39: WVR row 126 (anchor G=36210)
$ws.Range("H126").Value = 69668.62
$ws.Range("I126").Value = 74033
$ws.Range("K126").Value = 222099
$ws.Range("M126").Value = -219629

# hunk 40: WVR row 131 (anchor G=34723)
$ws.Range("H131").Value = 46773.6
$ws.Range("J131").Value = 46773.6
$ws.Range("L131").Value = 46773.6
$ws.Range("N131").Value = -56853.6

# hunk 41: WVR row 136 (anchor G=44031)
$ws.Range("H136").Value = 3961.6584
$ws.Range("I136").Value = 4105.9355
$ws.Range("J136").Value = 3514.4
$ws.Range("K136").Value = 12317.8065
$ws.Range("L136").Value = 10543.2
$ws.Range("M136").Value = -9767.806499999999
$ws.Range("N136").Value = -15643.2

